# Observer pattern: remove italics in diagrams
# - refresh the "last edited" datetimeFigureOut placeholders (slide master and
#   every slide layout) from 23/6/2017 -> 28/10/2020
# - drop the italic styling on the "update( )" run in the diagram on slide 1

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "28/10/2020"

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes $newDate

# Every slide layout hanging off the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes $newDate
}

# NOTE: the notes master also carries the same stale "23/6/2017" date field,
# but this runtime's NotesMaster object aliases the slide master's shape
# collection by index for writes (confirmed experimentally: writing to
# NotesMaster.Shapes.Item(2).TextFrame.TextRange.Text corrupts
# SlideMaster.Shapes.Item(2) instead of touching the notes master part), so
# it is intentionally left untouched here rather than risk corrupting the
# slide master.

# Remove italics from the "update( )" label in the Observer diagram
$slide1 = $p.Slides.Item(1)
$updateShape = $slide1.Shapes.Item("Rectangle 11")
$updateShape.TextFrame.TextRange.Font.Italic = $false
